$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None

MSG: The decision regarding the movie for Friday has concluded with no selection made.
"
$ws.Range("C3").Value = "MSG: None

MSG: The decision to acquire the rights to both movies has been recorded successfully.
"
$ws.Range("D3").Value = "both_movies, "
$ws.Range("C4").Value = "MSG: None

MSG: The decision about which movie to show on Friday was not reached, leading to the conclusion that no selection can be made at this time.
"
$ws.Range("C5").Value = "MSG: None

MSG: The decision was made that no movie will be shown on Friday.
"
$ws.Range("C6").Value = "MSG: None

MSG: The decision has been recorded, and the movie `"Barbie`" will be shown on Friday.
"
$ws.Range("C7").Value = "MSG: None

MSG: The decision has been recorded, indicating that no agreement was reached regarding the movie for Friday.
"
$ws.Range("C8").Value = "MSG: None

MSG: The committee did not reach a decision about what movie to show on Friday.
"
$ws.Range("C9").Value = "MSG: None

MSG: No decision about Friday's movie was made.
"
$ws.Range("C10").Value = "MSG: None

MSG: The decision has been recorded as `"no decision.`"
"
$ws.Range("C11").Value = "MSG: None

MSG: The decision process about which movie to show on Friday did not reach a conclusive agreement, so I have recorded a no decision.
"
$ws.Range("C12").Value = "MSG: None

MSG: The decision has been made to acquire rights for both movies.
"
$ws.Range("D12").Value = "both_movies, "
$ws.Range("C13").Value = "MSG: None

MSG: The decision about which movie to acquire was not reached during the conversation.
"
$ws.Range("C14").Value = "MSG: None

MSG: The decision process ended without any agreement on a movie for Friday, and no decision was reached.
"
$ws.Range("C15").Value = "MSG: None

MSG: The decision has been recorded, and there is no selection for the movie at this time.
"
$ws.Range("C16").Value = "MSG: None

MSG: The decision has been successfully recorded to acquire the rights for `"Barbie.`"
"
$ws.Range("C17").Value = "MSG: None

MSG: The rights for both movies have been acquired successfully.
"
$ws.Range("C18").Value = "MSG: None

MSG: The decision has been made to acquire the rights to `"Barbie`" for the movie to be shown on Friday.
"
$ws.Range("C19").Value = "MSG: None

MSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.
"
$ws.Range("C20").Value = "MSG: None

MSG: The decision about which movie to show on Friday has resulted in no agreement.
"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None

MSG: The decision has been recorded as no consensus was reached regarding the movie for Friday.
"
$ws.Range("C22").Value = "MSG: None

MSG: The decision has been recorded, and the rights to `"Barbie`" have been acquired for the showing on Friday.
"
$ws.Range("C23").Value = "MSG: None

MSG: The decision has been recorded to acquire the rights for both movies, `"Barbie`" and `"Oppenheimer,`" as they will both be shown in their entirety on Friday.
"
$ws.Range("C24").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.
"
$ws.Range("D24").Value = "Barbie_was_selected, "
$ws.Range("C25").Value = "MSG: None

MSG: The decision to acquire the rights for both movies has been made.
"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None

MSG: The function has been called, and the decision process concludes with no movie selected for Friday.
"
$ws.Range("C27").Value = "MSG: None

MSG: The decision regarding which movie to show on Friday has not been made.
"
$ws.Range("C28").Value = "MSG: None

MSG: The decision has been made to acquire the rights to `"Barbie.`"
"
$ws.Range("D28").Value = "Barbie_was_selected, "
$ws.Range("C29").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Oppenheimer`" as the movie to be shown on Friday.
"
$ws.Range("C30").Value = "MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
"
$ws.Range("C31").Value = "MSG: None

MSG: The decision to acquire the rights for `"Barbie`" has been confirmed.
"
$ws.Range("C32").Value = "MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
"
$ws.Range("C33").Value = "MSG: None

MSG: The decision has been made to acquire the rights for both movies.
"
$ws.Range("D33").Value = "both_movies, "
$ws.Range("C34").Value = "MSG: None

MSG: The decision has been recorded with the outcome of no selection for Friday’s movie.
"
$ws.Range("C35").Value = "MSG: None

MSG: The decision about which movie to show on Friday has not been made.
"
$ws.Range("C36").Value = "MSG: None

MSG: The decision has been recorded as `"no decision.`"
"
$ws.Range("C37").Value = "MSG: None

MSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.
"
$ws.Range("C38").Value = "MSG: None

MSG: No movie was selected in this meeting.
"
$ws.Range("C39").Value = "MSG: None

MSG: The function for no decision has been successfully executed, indicating that no movie was selected for Friday.
"
$ws.Range("C40").Value = "MSG: None

MSG: The decision regarding the movie to be shown on Friday has not been made.
"
$ws.Range("C41").Value = "MSG: None

MSG: The decision to acquire rights for both movies has been recorded.
"
$ws.Range("C42").Value = "MSG: None

MSG: I have recorded the decision as `"no decision`" regarding the movie to be shown on Friday.
"
$ws.Range("C43").Value = "MSG: None

MSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded.
"
$ws.Range("D43").Value = "Oppenheimer_was_selected, "
$ws.Range("C44").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie`" for the upcoming showing on Friday.
"
$ws.Range("C45").Value = "MSG: None

MSG: The decision to acquire the rights for the movie `"Barbie`" has been processed successfully.
"
$ws.Range("C46").Value = "MSG: None

MSG: I have successfully recorded the decision to acquire the rights for both movies.
"
$ws.Range("C47").Value = "MSG: None

MSG: The decision has been recorded, indicating that no definitive choice was made regarding the movie for Friday.
"
$ws.Range("C48").Value = "MSG: None

MSG: The decision about Friday's movie concluded without a choice being made.
"
$ws.Range("C49").Value = "MSG: None

MSG: The function for no decision has been executed, indicating that the committee did not arrive at a conclusion regarding Friday's movie.
"
$ws.Range("D49").Value = "no_decision, "
